# Refresh the activity report data (B2:L18) with the latest schedule-based
# timeframe values, as produced by the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(984,13,764.322,1112.327,66,12,3905,20215,34,262,6742),
    @(267,0,0,0,0,0,0,0,0,0,0),
    @(15,1,164,201,37,0,1841,1411,3,0,0),
    @(1098,2,177,192,12,0,2722,25952,5,0,0),
    @(1854,0,540,565,27,3,4587,68539,22,0,0),
    @(137,0,63,70,5,0,714,6349,0,0,0),
    @(805,10,1041,1186,58,2,727,13992,11,75,920),
    @(266,0,7,7,0,0,0,5714,0,0,0),
    @(288,6,930,976,38,6,2620,14538,11,0,0),
    @(1779,0,260,266,6,0,607,11436,0,0,0),
    @(1393,0,98,131,33,0,5377,6717,2,0,0),
    @(869,0,425,430,2,2,833,11190,12,1,417),
    @(268,0,15,15,0,0,0,6000,1,0,0),
    @(489,0,93,94,1,1,1250,30003,2,0,0),
    @(191,0,4,4,0,0,0,0,0,0,0),
    @(1812,0,50,52,2,0,769,18975,1,0,0),
    @(3653,30,650,1656,145,8,3746,22119,17,825,18065)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        # Column B is index 2
        $col = $c + 2
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
